$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.7304773333333333
$ws.Range("H2").Value = 2.191432
$ws.Range("I2").Value = 0.03163269997405359
$ws.Range("J2").Value = 0.03163269997405359
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.05601
$ws.Range("N2").Value = 0.16803
$ws.Range("O2").Value = 0.02710547761971223
$ws.Range("P2").Value = 0.02710547761971223
$ws.Range("Q2").Value = 0.04091403544
$ws.Range("R2").Value = 0.36822631896
$ws.Range("S2").Value = 0.0008574194411977813
$ws.Range("T2").Value = 0.0008574194411977814

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.7304773333333333
$ws.Range("H3").Value = 2.191432
$ws.Range("I3").Value = 0.03163269997405359
$ws.Range("J3").Value = 0.03163269997405359
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.864751
$ws.Range("N3").Value = 5.594253
$ws.Range("O3").Value = 0.902427539668559
$ws.Range("P3").Value = 0.9024275396685592
$ws.Range("Q3").Value = 1.362158337810667
$ws.Range("R3").Value = 12.259425040296
$ws.Range("S3").Value = 0.02854621961065888
$ws.Range("T3").Value = 0.02854621961065888

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.7304773333333333
$ws.Range("H4").Value = 2.191432
$ws.Range("I4").Value = 0.03163269997405359
$ws.Range("J4").Value = 0.03163269997405359
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.145611
$ws.Range("N4").Value = 0.436833
$ws.Range("O4").Value = 0.07046698271172858
$ws.Range("P4").Value = 0.07046698271172858
$ws.Range("Q4").Value = 0.106365534984
$ws.Range("R4").Value = 0.957289814856
$ws.Range("S4").Value = 0.002229060922196932
$ws.Range("T4").Value = 0.002229060922196932

$ws.Range("I5").Value = 0.4074771110502447
$ws.Range("J5").Value = 0.4074771110502448
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.05601
$ws.Range("N5").Value = 0.16803
$ws.Range("O5").Value = 0.02710547761971223
$ws.Range("P5").Value = 0.02710547761971223
$ws.Range("Q5").Value = 0.52703477655
$ws.Range("R5").Value = 4.743312988950001
$ws.Range("S5").Value = 0.0110448617141174
$ws.Range("T5").Value = 0.01104486171411741

$ws.Range("I6").Value = 0.4074771110502447
$ws.Range("J6").Value = 0.4074771110502448
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.864751
$ws.Range("N6").Value = 5.594253
$ws.Range("O6").Value = 0.902427539668559
$ws.Range("P6").Value = 0.9024275396685592
$ws.Range("S6").Value = 0.3677185667963245
$ws.Range("T6").Value = 0.3677185667963246

$ws.Range("I7").Value = 0.4074771110502447
$ws.Range("J7").Value = 0.4074771110502448
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.145611
$ws.Range("N7").Value = 0.436833
$ws.Range("O7").Value = 0.07046698271172858
$ws.Range("P7").Value = 0.07046698271172858
$ws.Range("Q7").Value = 1.370149274205
$ws.Range("R7").Value = 12.331343467845
$ws.Range("S7").Value = 0.0287136825398027
$ws.Range("T7").Value = 0.0287136825398027

$ws.Range("G8").Value = 12.95234266666667
$ws.Range("H8").Value = 38.857028
$ws.Range("I8").Value = 0.5608901889757016
$ws.Range("J8").Value = 0.5608901889757018
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.05601
$ws.Range("N8").Value = 0.16803
$ws.Range("O8").Value = 0.02710547761971223
$ws.Range("P8").Value = 0.02710547761971223
$ws.Range("Q8").Value = 0.72546071276
$ws.Range("R8").Value = 6.52914641484
$ws.Range("S8").Value = 0.01520319646439704
$ws.Range("T8").Value = 0.01520319646439705

$ws.Range("G9").Value = 12.95234266666667
$ws.Range("H9").Value = 38.857028
$ws.Range("I9").Value = 0.5608901889757016
$ws.Range("J9").Value = 0.5608901889757018
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.864751
$ws.Range("N9").Value = 5.594253
$ws.Range("O9").Value = 0.902427539668559
$ws.Range("P9").Value = 0.9024275396685592
$ws.Range("Q9").Value = 24.15289394000933
$ws.Range("R9").Value = 217.376045460084
$ws.Range("S9").Value = 0.5061627532615756
$ws.Range("T9").Value = 0.5061627532615758

$ws.Range("G10").Value = 12.95234266666667
$ws.Range("H10").Value = 38.857028
$ws.Range("I10").Value = 0.5608901889757016
$ws.Range("J10").Value = 0.5608901889757018
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.145611
$ws.Range("N10").Value = 0.436833
$ws.Range("O10").Value = 0.07046698271172858
$ws.Range("P10").Value = 0.07046698271172858
$ws.Range("Q10").Value = 1.886003568036
$ws.Range("R10").Value = 16.974032112324
$ws.Range("S10").Value = 0.03952423924972894
$ws.Range("T10").Value = 0.03952423924972895
